$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values for the Dkk2 -> Lrp6 ligand-receptor pair table.
# The cluster pairing now spans ECs and FAPs (sending) against
# ECs / FAPs / MuSCs (target), six rows total (previously three, FAPs only).

$data = @(
    @{ Row = 2;  A = "ECs";  B = "Dkk2"; C = "Lrp6"; D = "ECs";
       E = 1; F = 0.3333333333333333; G = 0.1293116666666667; H = 0.387935;
       I = 0.0151634279701504; J = 0.0151634279701504;
       K = 3; L = 1; M = 12.40685866666667; N = 37.220576;
       O = 0.1720325859617629; P = 0.1720325859617629;
       Q = 1.604351572284444; R = 14.43916415056;
       S = 0.0026086037257499; T = 0.0026086037257499 },

    @{ Row = 3;  A = "ECs";  B = "Dkk2"; C = "Lrp6"; D = "FAPs";
       E = 1; F = 0.3333333333333333; G = 0.1293116666666667; H = 0.387935;
       I = 0.0151634279701504; J = 0.0151634279701504;
       K = 3; L = 1; M = 43.97212233333332; N = 131.916367;
       O = 0.6097142007069145; P = 0.6097142007069145;
       Q = 5.686108425793887; R = 51.17497583214499;
       S = 0.009245357364797122; T = 0.009245357364797124 },

    @{ Row = 4;  A = "ECs";  B = "Dkk2"; C = "Lrp6"; D = "MuSCs";
       E = 1; F = 0.3333333333333333; G = 0.1293116666666667; H = 0.387935;
       I = 0.0151634279701504; J = 0.0151634279701504;
       K = 3; L = 1; M = 15.740255; N = 47.220765;
       O = 0.2182532133313226; P = 0.2182532133313226;
       Q = 2.035398607808333; R = 18.318587470275;
       S = 0.00330946687960338; T = 0.00330946687960338 },

    @{ Row = 5;  A = "FAPs"; B = "Dkk2"; C = "Lrp6"; D = "ECs";
       E = 3; F = 1; G = 8.398553333333334; H = 25.19566;
       I = 0.9848365720298496; J = 0.9848365720298496;
       K = 3; L = 1; M = 12.40685866666667; N = 37.220576;
       O = 0.1720325859617629; P = 0.1720325859617629;
       Q = 104.1996642111289; R = 937.7969779001601;
       S = 0.169423982236013; T = 0.169423982236013 },

    @{ Row = 6;  A = "FAPs"; B = "Dkk2"; C = "Lrp6"; D = "FAPs";
       E = 3; F = 1; G = 8.398553333333334; H = 25.19566;
       I = 0.9848365720298496; J = 0.9848365720298496;
       K = 3; L = 1; M = 43.97212233333332; N = 131.916367;
       O = 0.6097142007069145; P = 0.6097142007069145;
       Q = 369.3022145963577; R = 3323.719931367219;
       S = 0.6004688433421174; T = 0.6004688433421174 },

    @{ Row = 7;  A = "FAPs"; B = "Dkk2"; C = "Lrp6"; D = "MuSCs";
       E = 3; F = 1; G = 8.398553333333334; H = 25.19566;
       I = 0.9848365720298496; J = 0.9848365720298496;
       K = 3; L = 1; M = 15.740255; N = 47.220765;
       O = 0.2182532133313226; P = 0.2182532133313226;
       Q = 132.1953710977667; R = 1189.7583398799;
       S = 0.2149437464517192; T = 0.2149437464517192 }
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($rec in $data) {
    $r = $rec.Row
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value = $rec[$col]
    }
}
